$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990207675/1749682509529/zqlwK29jaz-camera-video-segment-driver-1749682514529.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSJKDNUWPN%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T180409Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJHMEUCIQCvl4NQn5TT%2FgwMKsU%2FrQI%2FrDRogR9DOzp7RiBXhb5QUgIgZi7NgF69wz7DSfRbYYfJYD4PO28LEvFmArdnUYjajq8q5gMI8v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDNGjs10ctG6giDFXlSq6A6vAEeEhpPIDemi0zEOC%2Bxsm3YILRZw2zy4rIZq7tO4WB8Rww%2Bwnp3OaGROw9f61kRZDHBCTssgjO%2Fe3dT26FunNPT4bbnvEpbtExa51%2Berv5N4QRXbVBxqTKHrE%2B5J%2FOoNqAerO72NQDyNWRVLAOJN40UWWGWjM5AOReNvrqKr7jo1Nv258pKoo1pvrcubtl7f0ED%2BNOxE2dvnh3IUAxNfclWCch%2FeKHaip5DolQ%2B2JfCBByHl%2BGyB6618yADBBACX8wqkGsJB2CEI329OdmZf1LMcg5ylIXEPFE5LOEhhaTPGlwYYqk4Jq%2BWMzURUmGHSZ1a7zAdvt0uMMYRXOpwwyJNsOhU4sD78h9GG48K4oJsKJJMdC8SOphuMbpwgguy3LodK%2FjAdFusnHzekyvIBf%2BxMZGvcvbE72ckb%2B6U0VaWTUaK4rAlwYDM8VGL%2Bi3HbmJ%2FpkJjLBQsurdryf394leOmdITb36yIoWFpq6lS4mZYQxAe9quHe8ly86%2B1LD9eAlGAwx5mGd7CktE0vxE6ZOL4fo3VKIUc73eGWPI9LcXAGLX3v7x6SQRG%2F%2BWXA2aub%2FIkfEOJQc18wtZGswgY6pQHQzLNPWiqpUc29UF%2B7OcSc5PlOrQut6caJcsyeGc%2Fo9Jx4LGEZiQxKYSUsPm5XVf%2B73YErGYh0cAZPWPOUnYdkasKDl%2B2h1GL722%2Br7nssa6QeQj12SGeCv4a57H93lMjwVEw1gvIxsDnmkuXS7GvD0AFMomlsr47TqRlAmpby%2BaRh9xGUw5T9wscBtn%2Fqll%2FAXtsxCW0EwAzeCnil8sTOI%2FeMlys%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2002%3A04%3A09%20GMT&X-Amz-Signature=c600b4a24955586e7f12d24095fd60868528ae9317001577a032a38f8b83382b'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990207675/1749682509529/JfEAwSUnRy-camera-video-segment-1749682514529.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSJKDNUWPN%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T180409Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJHMEUCIQCvl4NQn5TT%2FgwMKsU%2FrQI%2FrDRogR9DOzp7RiBXhb5QUgIgZi7NgF69wz7DSfRbYYfJYD4PO28LEvFmArdnUYjajq8q5gMI8v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDNGjs10ctG6giDFXlSq6A6vAEeEhpPIDemi0zEOC%2Bxsm3YILRZw2zy4rIZq7tO4WB8Rww%2Bwnp3OaGROw9f61kRZDHBCTssgjO%2Fe3dT26FunNPT4bbnvEpbtExa51%2Berv5N4QRXbVBxqTKHrE%2B5J%2FOoNqAerO72NQDyNWRVLAOJN40UWWGWjM5AOReNvrqKr7jo1Nv258pKoo1pvrcubtl7f0ED%2BNOxE2dvnh3IUAxNfclWCch%2FeKHaip5DolQ%2B2JfCBByHl%2BGyB6618yADBBACX8wqkGsJB2CEI329OdmZf1LMcg5ylIXEPFE5LOEhhaTPGlwYYqk4Jq%2BWMzURUmGHSZ1a7zAdvt0uMMYRXOpwwyJNsOhU4sD78h9GG48K4oJsKJJMdC8SOphuMbpwgguy3LodK%2FjAdFusnHzekyvIBf%2BxMZGvcvbE72ckb%2B6U0VaWTUaK4rAlwYDM8VGL%2Bi3HbmJ%2FpkJjLBQsurdryf394leOmdITb36yIoWFpq6lS4mZYQxAe9quHe8ly86%2B1LD9eAlGAwx5mGd7CktE0vxE6ZOL4fo3VKIUc73eGWPI9LcXAGLX3v7x6SQRG%2F%2BWXA2aub%2FIkfEOJQc18wtZGswgY6pQHQzLNPWiqpUc29UF%2B7OcSc5PlOrQut6caJcsyeGc%2Fo9Jx4LGEZiQxKYSUsPm5XVf%2B73YErGYh0cAZPWPOUnYdkasKDl%2B2h1GL722%2Br7nssa6QeQj12SGeCv4a57H93lMjwVEw1gvIxsDnmkuXS7GvD0AFMomlsr47TqRlAmpby%2BaRh9xGUw5T9wscBtn%2Fqll%2FAXtsxCW0EwAzeCnil8sTOI%2FeMlys%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2002%3A04%3A09%20GMT&X-Amz-Signature=16aff78db9b0e8999042c68dbe012bcd3f1db25ed07a0a183f2fda99aa56cc51'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991154589/1749681707753/UtIzP4FPML-camera-video-segment-driver-1749681712753.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSJKDNUWPN%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T180409Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJHMEUCIQCvl4NQn5TT%2FgwMKsU%2FrQI%2FrDRogR9DOzp7RiBXhb5QUgIgZi7NgF69wz7DSfRbYYfJYD4PO28LEvFmArdnUYjajq8q5gMI8v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDNGjs10ctG6giDFXlSq6A6vAEeEhpPIDemi0zEOC%2Bxsm3YILRZw2zy4rIZq7tO4WB8Rww%2Bwnp3OaGROw9f61kRZDHBCTssgjO%2Fe3dT26FunNPT4bbnvEpbtExa51%2Berv5N4QRXbVBxqTKHrE%2B5J%2FOoNqAerO72NQDyNWRVLAOJN40UWWGWjM5AOReNvrqKr7jo1Nv258pKoo1pvrcubtl7f0ED%2BNOxE2dvnh3IUAxNfclWCch%2FeKHaip5DolQ%2B2JfCBByHl%2BGyB6618yADBBACX8wqkGsJB2CEI329OdmZf1LMcg5ylIXEPFE5LOEhhaTPGlwYYqk4Jq%2BWMzURUmGHSZ1a7zAdvt0uMMYRXOpwwyJNsOhU4sD78h9GG48K4oJsKJJMdC8SOphuMbpwgguy3LodK%2FjAdFusnHzekyvIBf%2BxMZGvcvbE72ckb%2B6U0VaWTUaK4rAlwYDM8VGL%2Bi3HbmJ%2FpkJjLBQsurdryf394leOmdITb36yIoWFpq6lS4mZYQxAe9quHe8ly86%2B1LD9eAlGAwx5mGd7CktE0vxE6ZOL4fo3VKIUc73eGWPI9LcXAGLX3v7x6SQRG%2F%2BWXA2aub%2FIkfEOJQc18wtZGswgY6pQHQzLNPWiqpUc29UF%2B7OcSc5PlOrQut6caJcsyeGc%2Fo9Jx4LGEZiQxKYSUsPm5XVf%2B73YErGYh0cAZPWPOUnYdkasKDl%2B2h1GL722%2Br7nssa6QeQj12SGeCv4a57H93lMjwVEw1gvIxsDnmkuXS7GvD0AFMomlsr47TqRlAmpby%2BaRh9xGUw5T9wscBtn%2Fqll%2FAXtsxCW0EwAzeCnil8sTOI%2FeMlys%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2002%3A04%3A09%20GMT&X-Amz-Signature=ce8978277e81ded7b4357b6202f86bb331a662926bc45d6597e6d6692c1c321b'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991154589/1749681707753/418ikFdSWa-camera-video-segment-1749681712753.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSJKDNUWPN%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T180409Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJHMEUCIQCvl4NQn5TT%2FgwMKsU%2FrQI%2FrDRogR9DOzp7RiBXhb5QUgIgZi7NgF69wz7DSfRbYYfJYD4PO28LEvFmArdnUYjajq8q5gMI8v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDNGjs10ctG6giDFXlSq6A6vAEeEhpPIDemi0zEOC%2Bxsm3YILRZw2zy4rIZq7tO4WB8Rww%2Bwnp3OaGROw9f61kRZDHBCTssgjO%2Fe3dT26FunNPT4bbnvEpbtExa51%2Berv5N4QRXbVBxqTKHrE%2B5J%2FOoNqAerO72NQDyNWRVLAOJN40UWWGWjM5AOReNvrqKr7jo1Nv258pKoo1pvrcubtl7f0ED%2BNOxE2dvnh3IUAxNfclWCch%2FeKHaip5DolQ%2B2JfCBByHl%2BGyB6618yADBBACX8wqkGsJB2CEI329OdmZf1LMcg5ylIXEPFE5LOEhhaTPGlwYYqk4Jq%2BWMzURUmGHSZ1a7zAdvt0uMMYRXOpwwyJNsOhU4sD78h9GG48K4oJsKJJMdC8SOphuMbpwgguy3LodK%2FjAdFusnHzekyvIBf%2BxMZGvcvbE72ckb%2B6U0VaWTUaK4rAlwYDM8VGL%2Bi3HbmJ%2FpkJjLBQsurdryf394leOmdITb36yIoWFpq6lS4mZYQxAe9quHe8ly86%2B1LD9eAlGAwx5mGd7CktE0vxE6ZOL4fo3VKIUc73eGWPI9LcXAGLX3v7x6SQRG%2F%2BWXA2aub%2FIkfEOJQc18wtZGswgY6pQHQzLNPWiqpUc29UF%2B7OcSc5PlOrQut6caJcsyeGc%2Fo9Jx4LGEZiQxKYSUsPm5XVf%2B73YErGYh0cAZPWPOUnYdkasKDl%2B2h1GL722%2Br7nssa6QeQj12SGeCv4a57H93lMjwVEw1gvIxsDnmkuXS7GvD0AFMomlsr47TqRlAmpby%2BaRh9xGUw5T9wscBtn%2Fqll%2FAXtsxCW0EwAzeCnil8sTOI%2FeMlys%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2002%3A04%3A09%20GMT&X-Amz-Signature=759cb0c7c84d9f7074426229b803b13691077754be1cb42158c3ca5f5b0a0dba'
